# Update the "Lottery Results" sheet with the latest draw data.
# Rows 2-8 get new winning numbers/prizes, and rows 9-13 are brand new
# entries appended below the existing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Ticket Number" (column B) must stay text even though the values look
# numeric - force a Text number format before writing so Excel doesn't
# silently convert it to a number (column C keeps the numeric duplicate).
$ws.Range("B2:B13").NumberFormat = "@"

# --- Row 2 ---
$ws.Range("B2").Value = "12024"
$ws.Range("C2").Value = 12024
$ws.Range("E2").Value = "Futura Pressure Cooker"
$ws.Range("F2").Value = "/static/prizes/futura_pressure_cooker.jpg"

# --- Row 3 ---
$ws.Range("B3").Value = "17227"
$ws.Range("C3").Value = 17227
$ws.Range("D3").Value = "Lumbini - Bhairahawa"
$ws.Range("E3").Value = "Futura Pressure Cooker"
$ws.Range("F3").Value = "/static/prizes/futura_pressure_cooker.jpg"

# --- Row 4 ---
$ws.Range("B4").Value = "11586"
$ws.Range("C4").Value = 11586
$ws.Range("D4").Value = "Bagmati"

# --- Row 5 ---
$ws.Range("B5").Value = "13407"
$ws.Range("C5").Value = 13407
$ws.Range("D5").Value = "Dang"

# --- Row 6 ---
$ws.Range("B6").Value = "14451"
$ws.Range("C6").Value = 14451
$ws.Range("D6").Value = "Gandaki"
$ws.Range("E6").Value = "Vacuum Cleaner"
$ws.Range("F6").Value = "/static/prizes/vacuum_cleaner.jpg"

# --- Row 7 ---
$ws.Range("B7").Value = "16472"
$ws.Range("C7").Value = 16472
$ws.Range("D7").Value = "Gandaki"

# --- Row 8 ---
$ws.Range("B8").Value = "11204"
$ws.Range("C8").Value = 11204
$ws.Range("D8").Value = "Bagmati"
$ws.Range("E8").Value = "Vacuum Cleaner"
$ws.Range("F8").Value = "/static/prizes/vacuum_cleaner.jpg"

# --- Row 9 (new) ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "10343"
$ws.Range("C9").Value = 10343
$ws.Range("D9").Value = "Mu Ka"
$ws.Range("E9").Value = "Futura Pressure Cooker"
$ws.Range("F9").Value = "/static/prizes/futura_pressure_cooker.jpg"

# --- Row 10 (new) ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "13590"
$ws.Range("C10").Value = 13590
$ws.Range("D10").Value = "Dang"
$ws.Range("E10").Value = "Vacuum Cleaner"
$ws.Range("F10").Value = "/static/prizes/vacuum_cleaner.jpg"

# --- Row 11 (new) ---
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "18428"
$ws.Range("C11").Value = 18428
$ws.Range("D11").Value = "Janakpur"
$ws.Range("E11").Value = "Samsung Washing Machine"
$ws.Range("F11").Value = "/static/prizes/samsung_washing_machine.jpg"

# --- Row 12 (new) ---
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "14280"
$ws.Range("C12").Value = 14280
$ws.Range("D12").Value = "Gandaki"
$ws.Range("E12").Value = "Futura Pressure Cooker"
$ws.Range("F12").Value = "/static/prizes/futura_pressure_cooker.jpg"

# --- Row 13 (new) ---
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "17124"
$ws.Range("C13").Value = 17124
$ws.Range("D13").Value = "Bagmati"
$ws.Range("E13").Value = "Futura Pressure Cooker"
$ws.Range("F13").Value = "/static/prizes/futura_pressure_cooker.jpg"
